# "Ready for test again"
# - Duplicate "Sheet 5" twice, appended after "Sheet 6", as "Sheet 7" and "Sheet 8".
# - Make "Sheet 4" the active sheet, scrolled/selected to D2, zoomed to 75%.
# - (As a side effect of switching the active sheet, "Sheet 3" loses its
#   tabSelected flag and its old scroll position.)

$wb = $excel.ActiveWorkbook

$sheet5 = $wb.Worksheets.Item("Sheet 5")
$sheet6 = $wb.Worksheets.Item("Sheet 6")

# First duplicate -> placed right after "Sheet 6", becomes "Sheet 7".
$sheet5.Copy($null, $sheet6)
$newSheet1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet1.Name = "Sheet 7"

# Second duplicate -> placed right after the first new sheet, becomes "Sheet 8".
$sheet5.Copy($null, $newSheet1)
$newSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet2.Name = "Sheet 8"

# Switch the active/visible sheet to "Sheet 4" and update its view.
$ws4 = $wb.Worksheets.Item("Sheet 4")
$ws4.Activate() | Out-Null
$ws4.Range("D2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75
